# Applies the "Add files via upload" edit: replaces the sample AP rows
# with a single generic placeholder row, removing rows 3 and 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with generic placeholder values.
$ws.Range("B2").Value = "00:00:00:00:00:00"
$ws.Range("C2").Value = "192.168.0.111"

# Remove the now-unused rows 3 and 4 entirely (shrinks used range to A1:C2).
$ws.Range("A3:C4").Delete()

# Match the saved selection from the edited file.
$ws.Range("A2").Select()
